$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 26461
$ws.Cells.Item(2, 5).Value = 512841095632
$ws.Cells.Item(2, 6).Value = 15295281760
$ws.Cells.Item(2, 7).Value = 0.36016

$ws.Cells.Item(3, 4).Value = 1808.22
$ws.Cells.Item(3, 5).Value = 217310665995
$ws.Cells.Item(3, 6).Value = 7560810565
$ws.Cells.Item(3, 7).Value = 0.39102

$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(4, 5).Value = 83048716680
$ws.Cells.Item(4, 6).Value = 16986580668
$ws.Cells.Item(4, 7).Value = 0.01672

$ws.Cells.Item(5, 4).Value = 304.59
$ws.Cells.Item(5, 5).Value = 48093350164
$ws.Cells.Item(5, 6).Value = 552148251
$ws.Cells.Item(5, 7).Value = -0.34153

$ws.Cells.Item(6, 4).Value = 0.999893
$ws.Cells.Item(6, 5).Value = 29046812679
$ws.Cells.Item(6, 6).Value = 4146045655
$ws.Cells.Item(6, 7).Value = 0.00342

$ws.Cells.Item(7, 4).Value = 0.453827
$ws.Cells.Item(7, 5).Value = 23510158592
$ws.Cells.Item(7, 6).Value = 823922217
$ws.Cells.Item(7, 7).Value = -0.0319

$ws.Cells.Item(8, 4).Value = 0.358416
$ws.Cells.Item(8, 5).Value = 12548450887
$ws.Cells.Item(8, 6).Value = 173579966
$ws.Cells.Item(8, 7).Value = -1.6919

$ws.Cells.Item(9, 4).Value = 1804.89
$ws.Cells.Item(9, 5).Value = 12079944766
$ws.Cells.Item(9, 6).Value = 7261597
$ws.Cells.Item(9, 7).Value = 0.25163

$ws.Cells.Item(10, 4).Value = 0.07091
$ws.Cells.Item(10, 5).Value = 9884459453
$ws.Cells.Item(10, 6).Value = 277203003
$ws.Cells.Item(10, 7).Value = -0.04145

$ws.Cells.Item(11, 4).Value = 0.888668
$ws.Cells.Item(11, 5).Value = 8241486339
$ws.Cells.Item(11, 6).Value = 293209644
$ws.Cells.Item(11, 7).Value = 1.73392

$ws.Cells.Item(12, 4).Value = 19.24
$ws.Cells.Item(12, 5).Value = 7618113002
$ws.Cells.Item(12, 6).Value = 235056524
$ws.Cells.Item(12, 7).Value = -0.08337

$ws.Cells.Item(13, 4).Value = 0.076998
$ws.Cells.Item(13, 5).Value = 6962525507
$ws.Cells.Item(13, 6).Value = 363493525
$ws.Cells.Item(13, 7).Value = -0.27074

$ws.Cells.Item(14, 4).Value = 5.25
$ws.Cells.Item(14, 5).Value = 6481987436
$ws.Cells.Item(14, 6).Value = 108728990
$ws.Cells.Item(14, 7).Value = -0.2123

$ws.Cells.Item(15, 4).Value = 85.62
$ws.Cells.Item(15, 5).Value = 6242787836
$ws.Cells.Item(15, 6).Value = 630622046
$ws.Cells.Item(15, 7).Value = -0.24829

$ws.Cells.Item(16, 4).Value = 0.999919
$ws.Cells.Item(16, 5).Value = 5301038364
$ws.Cells.Item(16, 6).Value = 1656697773
$ws.Cells.Item(16, 7).Value = -0.04618

$ws.Cells.Item(17, 5).Value = 5004629257
$ws.Cells.Item(17, 6).Value = 131936290
$ws.Cells.Item(17, 7).Value = -0.58525

$ws.Cells.Item(18, 4).Value = 14.08
$ws.Cells.Item(18, 5).Value = 4709807307
$ws.Cells.Item(18, 6).Value = 148626546
$ws.Cells.Item(18, 7).Value = -0.38017

$ws.Cells.Item(19, 4).Value = 1
$ws.Cells.Item(19, 5).Value = 4621598349
$ws.Cells.Item(19, 6).Value = 93803370
$ws.Cells.Item(19, 7).Value = -0.01675

$ws.Cells.Item(20, 4).Value = 26525
$ws.Cells.Item(20, 5).Value = 4130404968
$ws.Cells.Item(20, 6).Value = 81162016
$ws.Cells.Item(20, 7).Value = 0.43469

$ws.Cells.Item(21, 4).Value = 4.93
$ws.Cells.Item(21, 5).Value = 3714508265
$ws.Cells.Item(21, 6).Value = 35825947
$ws.Cells.Item(21, 7).Value = -0.81306

$ws.Cells.Item(22, 4).Value = 3.55
$ws.Cells.Item(22, 5).Value = 3312230830
$ws.Cells.Item(22, 6).Value = 211734
$ws.Cells.Item(22, 7).Value = 0.1712

$ws.Cells.Item(23, 4).Value = 6.28
$ws.Cells.Item(23, 5).Value = 3245978583
$ws.Cells.Item(23, 6).Value = 121588303
$ws.Cells.Item(23, 7).Value = -0.77695

$ws.Cells.Item(24, 4).Value = 10.47
$ws.Cells.Item(24, 5).Value = 3059269614
$ws.Cells.Item(24, 6).Value = 73757240
$ws.Cells.Item(24, 7).Value = 0.68104

$ws.Cells.Item(25, 4).Value = 1.92
$ws.Cells.Item(25, 5).Value = 2834567330
$ws.Cells.Item(25, 6).Value = 9350678
$ws.Cells.Item(25, 7).Value = -2.45678

$ws.Cells.Item(26, 4).Value = 46.21
$ws.Cells.Item(26, 5).Value = 2772979836
$ws.Cells.Item(26, 6).Value = 16865532
$ws.Cells.Item(26, 7).Value = 3.28009

$ws.Cells.Item(27, 4).Value = 151.49
$ws.Cells.Item(27, 5).Value = 2749238679
$ws.Cells.Item(27, 6).Value = 77370081
$ws.Cells.Item(27, 7).Value = 0.52234

$ws.Cells.Item(28, 4).Value = 17.71
$ws.Cells.Item(28, 5).Value = 2498819995
$ws.Cells.Item(28, 6).Value = 77224810
$ws.Cells.Item(28, 7).Value = -0.57694

$ws.Cells.Item(29, 4).Value = 0.086969
$ws.Cells.Item(29, 5).Value = 2329514531
$ws.Cells.Item(29, 6).Value = 32556750
$ws.Cells.Item(29, 7).Value = 0.54269

$ws.Cells.Item(30, 4).Value = 111.96
$ws.Cells.Item(30, 5).Value = 2170058440
$ws.Cells.Item(30, 6).Value = 57106195
$ws.Cells.Item(30, 7).Value = -0.20418

$ws.Cells.Item(31, 4).Value = 4.8
$ws.Cells.Item(31, 5).Value = 2096632017
$ws.Cells.Item(31, 6).Value = 23885608
$ws.Cells.Item(31, 7).Value = -0.5946

$ws.Cells.Item(32, 4).Value = 0.999465
$ws.Cells.Item(32, 5).Value = 2039836347
$ws.Cells.Item(32, 6).Value = 201366363
$ws.Cells.Item(32, 7).Value = -0.03849

$ws.Cells.Item(33, 5).Value = 1881920092
$ws.Cells.Item(33, 6).Value = 87034526
$ws.Cells.Item(33, 7).Value = -1.33831

$ws.Cells.Item(34, 5).Value = 1772048288
$ws.Cells.Item(34, 6).Value = 59205327
$ws.Cells.Item(34, 7).Value = -1.63658

$ws.Cells.Item(35, 4).Value = 0.050879
$ws.Cells.Item(35, 5).Value = 1599809380
$ws.Cells.Item(35, 6).Value = 16274494
$ws.Cells.Item(35, 7).Value = -0.375

$ws.Cells.Item(36, 4).Value = 8
$ws.Cells.Item(36, 5).Value = 1585958487
$ws.Cells.Item(36, 6).Value = 40217056
$ws.Cells.Item(36, 7).Value = -1.38525

$ws.Cells.Item(37, 4).Value = 0.059887
$ws.Cells.Item(37, 5).Value = 1512482196
$ws.Cells.Item(37, 6).Value = 5287148
$ws.Cells.Item(37, 7).Value = -0.28133

$ws.Cells.Item(38, 4).Value = 98.53
$ws.Cells.Item(38, 5).Value = 1435538827
$ws.Cells.Item(38, 6).Value = 17980109
$ws.Cells.Item(38, 7).Value = -1.7878

$ws.Cells.Item(39, 4).Value = 1.56
$ws.Cells.Item(39, 5).Value = 1423951244
$ws.Cells.Item(39, 6).Value = 58596041
$ws.Cells.Item(39, 7).Value = -1.46816

$ws.Cells.Item(40, 4).Value = 1.11
$ws.Cells.Item(40, 5).Value = 1410811743
$ws.Cells.Item(40, 6).Value = 168047151
$ws.Cells.Item(40, 7).Value = -0.66332

$ws.Cells.Item(41, 4).Value = 0.01929719
$ws.Cells.Item(41, 5).Value = 1402060808
$ws.Cells.Item(41, 6).Value = 38562861
$ws.Cells.Item(41, 7).Value = -0.49899

$ws.Cells.Item(42, 5).Value = 1191675051
$ws.Cells.Item(42, 6).Value = 64886397
$ws.Cells.Item(42, 7).Value = -3.33958

$ws.Cells.Item(43, 4).Value = 0.150811
$ws.Cells.Item(43, 5).Value = 1092287208
$ws.Cells.Item(43, 6).Value = 65603081
$ws.Cells.Item(43, 7).Value = -3.19802

$ws.Cells.Item(44, 2).Value = "GRT"
$ws.Cells.Item(44, 3).Value = "The Graph"
$ws.Cells.Item(44, 4).Value = 0.116712
$ws.Cells.Item(44, 5).Value = 1048598196
$ws.Cells.Item(44, 6).Value = 31429935
$ws.Cells.Item(44, 7).Value = 0.18855

$ws.Cells.Item(45, 2).Value = "GGTKN"
$ws.Cells.Item(45, 3).Value = "GGTKN"
$ws.Cells.Item(45, 4).Value = 0.092044
$ws.Cells.Item(45, 5).Value = 1046827992
$ws.Cells.Item(45, 6).Value = 68601
$ws.Cells.Item(45, 7).Value = 1.74621

$ws.Cells.Item(46, 4).Value = 1
$ws.Cells.Item(46, 5).Value = 1018951360
$ws.Cells.Item(46, 6).Value = 28000680
$ws.Cells.Item(46, 7).Value = 0.08774999999999999

$ws.Cells.Item(47, 4).Value = 0.9989980000000001
$ws.Cells.Item(47, 5).Value = 1003047199
$ws.Cells.Item(47, 6).Value = 15735583
$ws.Cells.Item(47, 7).Value = 0.0129

$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = "EDGT"
$ws.Cells.Item(48, 3).Value = "Edgecoin"
$ws.Cells.Item(48, 4).Value = 0.999876
$ws.Cells.Item(48, 5).Value = 999907301
$ws.Cells.Item(48, 6).Value = 5367459
$ws.Cells.Item(48, 7).Value = -0.01041

$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = "RNDR"
$ws.Cells.Item(49, 3).Value = "Render"
$ws.Cells.Item(49, 4).Value = 2.71
$ws.Cells.Item(49, 5).Value = 991826838
$ws.Cells.Item(49, 6).Value = 223721046
$ws.Cells.Item(49, 7).Value = 1.8508

$ws.Cells.Item(50, 4).Value = 0.509141
$ws.Cells.Item(50, 5).Value = 944138838
$ws.Cells.Item(50, 6).Value = 174435070
$ws.Cells.Item(50, 7).Value = 1.97515

$ws.Cells.Item(51, 4).Value = 0.84457
$ws.Cells.Item(51, 5).Value = 932574554
$ws.Cells.Item(51, 6).Value = 84661316
$ws.Cells.Item(51, 7).Value = 0.59019
